# edit.ps1
# Applies the changes described by the commit:
#   "removing modulo at quantizer function, changing mse per sample to
#    nromalized mse"
#
#  1. Refresh the cached "datetimeFigureOut" date field text (slide
#     master + all 11 slide layouts) from 2016-11-22 to 2016-12-06.
#  2. On slide 7 ("continue"), rewrite the two body paragraphs so that
#     every occurrence of "covariance between " becomes
#     "conditional variance between ", using an in-place
#     character-range replace (mirrors how a user would select just
#     that word in PowerPoint and retype it, which is what produces
#     the run-split seen in the target deck).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Date placeholder ("datetimeFigureOut" field) on master + layouts
# ---------------------------------------------------------------
$oldDate = "2016-11-22"
$newDate = "2016-12-06"

$m = $p.SlideMaster

# Master's own Date Placeholder.
for ($shi = 1; $shi -le $m.Shapes.Count; $shi++) {
    $sh = $m.Shapes.Item($shi)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

# Every custom (slide) layout also carries its own copy of the field.
for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $layout = $m.CustomLayouts.Item($li)
    for ($shi = 1; $shi -le $layout.Shapes.Count; $shi++) {
        $sh = $layout.Shapes.Item($shi)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------
# 2) Slide 7 body text: "covariance" -> "conditional variance"
# ---------------------------------------------------------------
$slide7 = $p.Slides.Item(7)
for ($shi = 1; $shi -le $slide7.Shapes.Count; $shi++) {
    $sh = $slide7.Shapes.Item($shi)
    if (-not $sh.HasTextFrame) { continue }
    $text = $sh.TextFrame.TextRange.Text
    if ($text -notlike "*covariance between*") { continue }

    $tr2 = $sh.TextFrame2.TextRange
    $old = "covariance between "
    $new = "conditional variance between "

    # Locate every occurrence in the ORIGINAL text first so the
    # offsets are stable, then replay the edits left-to-right,
    # shifting later offsets by the cumulative growth introduced by
    # earlier (already-applied) replacements.
    $orig = $tr2.Text
    $offsets = New-Object System.Collections.Generic.List[int]
    $searchFrom = 0
    while ($true) {
        $idx = $orig.IndexOf($old, $searchFrom)
        if ($idx -lt 0) { break }
        [void]$offsets.Add($idx)
        $searchFrom = $idx + $old.Length
    }

    $growth = $new.Length - $old.Length
    for ($k = 0; $k -lt $offsets.Count; $k++) {
        $start1based = $offsets[$k] + 1 + ($growth * $k)
        $range = $tr2.Characters($start1based, $old.Length)
        $range.Text = $new
    }
}
